$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 3022.5
$ws.Range("I18").Value = 696.6667
$ws.Range("K18").Value = 696.6667
$ws.Range("M18").Value = -412.6667
$ws.Range("H43").Value = 851.2222
$ws.Range("I43").Value = 864.75
$ws.Range("J43").Value = 840.4
$ws.Range("K43").Value = 864.75
$ws.Range("L43").Value = 840.4
$ws.Range("M43").Value = -795.75
$ws.Range("N43").Value = -978.4
$ws.Range("H64").Value = 6473.778
$ws.Range("I64").Value = 3600
$ws.Range("K64").Value = 3600
$ws.Range("M64").Value = -3352
$ws.Range("H67").Value = 6473.778
$ws.Range("I67").Value = 3600
$ws.Range("K67").Value = 3600
$ws.Range("M67").Value = -2742
$ws.Range("H70").Value = 4683.2
$ws.Range("I70").Value = 1875
$ws.Range("J70").Value = 5115.231
$ws.Range("K70").Value = 5625
$ws.Range("L70").Value = 15345.693
$ws.Range("M70").Value = -5355
$ws.Range("N70").Value = -15885.693
$ws.Range("H73").Value = 4683.2
$ws.Range("I73").Value = 1875
$ws.Range("J73").Value = 5115.231
$ws.Range("K73").Value = 5625
$ws.Range("L73").Value = 15345.693
$ws.Range("M73").Value = -4689
$ws.Range("N73").Value = -17217.693
$ws.Range("H76").Value = 4031.7144
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 4031.7144
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H98").Value = 1207.591
$ws.Range("I98").Value = 1108.7894
$ws.Range("J98").Value = 1833.3334
$ws.Range("K98").Value = 1108.7894
$ws.Range("L98").Value = 1833.3334
$ws.Range("M98").Value = 389.2106000000001
$ws.Range("N98").Value = -4829.3334
$ws.Range("H107").Value = 690.34485
$ws.Range("I107").Value = 657.75
$ws.Range("J107").Value = 846.8
$ws.Range("K107").Value = 657.75
$ws.Range("L107").Value = 846.8
$ws.Range("M107").Value = 1262.25
$ws.Range("N107").Value = -4686.8
$ws.Range("H112").Value = 73116.28999999999
$ws.Range("J112").Value = 92916.73
$ws.Range("L112").Value = 278750.19
$ws.Range("N112").Value = -280966.19
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H122").Value = 1207.591
$ws.Range("I122").Value = 1108.7894
$ws.Range("J122").Value = 1833.3334
$ws.Range("K122").Value = 3326.3682
$ws.Range("L122").Value = 5500.0002
$ws.Range("M122").Value = -876.3681999999999
$ws.Range("N122").Value = -10400.0002
$ws.Range("H132").Value = 6061578
$ws.Range("I132").Value = 7093194
$ws.Range("K132").Value = 21279582
$ws.Range("M132").Value = -21277052
$ws.Range("H135").Value = 963.84375
$ws.Range("I135").Value = 618.04346
$ws.Range("J135").Value = 1847.5555
$ws.Range("K135").Value = 5562.39114
$ws.Range("L135").Value = 16627.9995
$ws.Range("M135").Value = -3027.39114
$ws.Range("N135").Value = -21697.9995
$ws.Range("H137").Value = 1800.3572
$ws.Range("I137").Value = 1631.0435
$ws.Range("K137").Value = 4893.1305
$ws.Range("M137").Value = -2343.1305
$ws.Range("H141").Value = 2338.0625
$ws.Range("I141").Value = 1827.2667
$ws.Range("K141").Value = 5481.800099999999
$ws.Range("M141").Value = -301.8000999999995
# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2720.5774
$ws.Range("I32").Value = 2048.8462
$ws.Range("K32").Value = 2048.8462
$ws.Range("M32").Value = -1761.8462
$ws.Range("H45").Value = 5002.108
$ws.Range("I45").Value = 5002.2285
$ws.Range("K45").Value = 5002.2285
$ws.Range("M45").Value = -4625.2285
$ws.Range("H61").Value = 5738.75
$ws.Range("I61").Value = 5325.625
$ws.Range("J61").Value = 6151.875
$ws.Range("K61").Value = 5325.625
$ws.Range("L61").Value = 6151.875
$ws.Range("M61").Value = -5113.625
$ws.Range("N61").Value = -6575.875
$ws.Range("H74").Value = 5801.3
$ws.Range("I74").Value = 1254.7576
$ws.Range("J74").Value = 27235
$ws.Range("K74").Value = 1254.7576
$ws.Range("L74").Value = 27235
$ws.Range("M74").Value = -380.7575999999999
$ws.Range("N74").Value = -28983
$ws.Range("H77").Value = 5801.3
$ws.Range("I77").Value = 1254.7576
$ws.Range("J77").Value = 27235
$ws.Range("K77").Value = 6273.788
$ws.Range("L77").Value = 136175
$ws.Range("M77").Value = -1905.788
$ws.Range("N77").Value = -144911
$ws.Range("H132").Value = 3978.0625
$ws.Range("I132").Value = 2180.8
$ws.Range("J132").Value = 6973.5
$ws.Range("K132").Value = 6542.400000000001
$ws.Range("L132").Value = 20920.5
$ws.Range("M132").Value = -4012.400000000001
$ws.Range("N132").Value = -25980.5
$ws.Range("H136").Value = 5738.75
$ws.Range("I136").Value = 5325.625
$ws.Range("J136").Value = 6151.875
$ws.Range("K136").Value = 15976.875
$ws.Range("L136").Value = 18455.625
$ws.Range("M136").Value = -13426.875
$ws.Range("N136").Value = -23555.625
# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H111").Value = 45000
$ws.Range("J111").Value = 45000
$ws.Range("L111").Value = 45000
$ws.Range("N111").Value = -53180
$ws.Range("H134").Value = 2129.7666
$ws.Range("I134").Value = 2134.2415
$ws.Range("K134").Value = 6402.7245
$ws.Range("M134").Value = -3867.7245
# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 3401.6667
$ws.Range("I2").Value = 100
$ws.Range("K2").Value = 100
$ws.Range("M2").Value = 13
# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 8729.333000000001
$ws.Range("J5").Value = 10980.714
$ws.Range("L5").Value = 32942.142
$ws.Range("N5").Value = -33166.142
$ws.Range("H68").Value = 4548216.5
$ws.Range("I68").Value = 2049.75
$ws.Range("J68").Value = 7146025.5
$ws.Range("K68").Value = 6149.25
$ws.Range("L68").Value = 21438076.5
$ws.Range("M68").Value = -5338.25
$ws.Range("N68").Value = -21439698.5
$ws.Range("H71").Value = 4548216.5
$ws.Range("I71").Value = 2049.75
$ws.Range("J71").Value = 7146025.5
$ws.Range("K71").Value = 18447.75
$ws.Range("L71").Value = 64314229.5
$ws.Range("M71").Value = -14391.75
$ws.Range("N71").Value = -64322341.5
$ws.Range("H135").Value = 8729.333000000001
$ws.Range("J135").Value = 10980.714
$ws.Range("L135").Value = 98826.42600000001
$ws.Range("N135").Value = -103896.426
# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 5000
$ws.Range("I4").Value = 5000
$ws.Range("K4").Value = 5000
$ws.Range("M4").Value = -4888
$ws.Range("H126").Value = 13319.774
$ws.Range("I126").Value = 16535.348
$ws.Range("J126").Value = 4075
$ws.Range("K126").Value = 49606.04400000001
$ws.Range("L126").Value = 12225
$ws.Range("M126").Value = -47136.04400000001
$ws.Range("N126").Value = -17165
# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 3166.6667
$ws.Range("I4").Value = 3000
$ws.Range("K4").Value = 3000
$ws.Range("M4").Value = -2887
$ws.Range("H16").Value = 20836056
$ws.Range("I16").Value = 25001932
$ws.Range("J16").Value = 6674.75
$ws.Range("K16").Value = 25001932
$ws.Range("L16").Value = 6674.75
$ws.Range("M16").Value = -25001762
$ws.Range("N16").Value = -7014.75
$ws.Range("H28").Value = 3166.6667
$ws.Range("I28").Value = 3000
$ws.Range("K28").Value = 3000
$ws.Range("M28").Value = -2768
$ws.Range("H37").Value = 3166.6667
$ws.Range("I37").Value = 3000
$ws.Range("K37").Value = 3000
$ws.Range("M37").Value = -2893
$ws.Range("H46").Value = 2310.6
$ws.Range("I46").Value = 1229.8334
$ws.Range("J46").Value = 3031.111
$ws.Range("K46").Value = 1229.8334
$ws.Range("L46").Value = 3031.111
$ws.Range("M46").Value = -1041.8334
$ws.Range("N46").Value = -3407.111
$ws.Range("H55").Value = 614.8461
$ws.Range("I55").Value = 628.2222
$ws.Range("K55").Value = 628.2222
$ws.Range("M55").Value = -455.2222
$ws.Range("H132").Value = 3833.742
$ws.Range("I132").Value = 3414.64
$ws.Range("K132").Value = 10243.92
$ws.Range("M132").Value = -7713.92
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 1714.7273
$ws.Range("I6").Value = 425
$ws.Range("J6").Value = 2198.375
$ws.Range("K6").Value = 425
$ws.Range("L6").Value = 2198.375
$ws.Range("M6").Value = -310
$ws.Range("N6").Value = -2428.375
$ws.Range("H14").Value = 634.7083
$ws.Range("I14").Value = 606.381
$ws.Range("J14").Value = 833
$ws.Range("K14").Value = 606.381
$ws.Range("L14").Value = 833
$ws.Range("M14").Value = -438.381
$ws.Range("N14").Value = -1169
$ws.Range("H113").Value = 904.6875
$ws.Range("I113").Value = 817.6
$ws.Range("J113").Value = 944.2727
$ws.Range("K113").Value = 2452.8
$ws.Range("L113").Value = 2832.8181
$ws.Range("M113").Value = -282.8000000000002
$ws.Range("N113").Value = -7172.8181
$ws.Range("H141").Value = 200357.5
$ws.Range("J141").Value = 200357.5
$ws.Range("L141").Value = 200357.5
$ws.Range("N141").Value = -210717.5
